# Make tweet query properties consistent with ontology properties
# (rename relation labels used in the D/E/F "predecessor atom / relation / successor atom"
#  columns so they match the renamed ontology properties, and propagate the
#  "causes" -> "causesCondition" rewording into the affected tweet texts.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$apos = [char]0x2019   # RIGHT SINGLE QUOTATION MARK (Alzheimer's)
$bull = [char]0x2022   # BULLET

# --- relation-column (E) and a couple of successor-column (F) renames ---
$ws.Range("E2").Value  = "causesCondition"
$ws.Range("E3").Value  = "preventsCondition"
$ws.Range("E4").Value  = "causesCondition"
$ws.Range("E5").Value  = "not sportCausesInjury"
$ws.Range("E6").Value  = "not sportCausesInjury"
$ws.Range("E7").Value  = "causesCondition"
$ws.Range("E8").Value  = "causesCondition"
$ws.Range("E9").Value  = "CausesNutrientState"
$ws.Range("F9").Value  = "HealthyNutritionState"
$ws.Range("E10").Value = "CausesNutrientState"
$ws.Range("F10").Value = "HealthyNutritionState"
$ws.Range("E12").Value = "Contains"
$ws.Range("E13").Value = "preventsCondition"
$ws.Range("E14").Value = "preventsCondition"
$ws.Range("E15").Value = "preventsCondition"
$ws.Range("E16").Value = "causesCondition"
$ws.Range("E17").Value = "causesCondition"

# --- tweet text bodies: "causes" -> "causesCondition" ---
$bigFatSurpriseTweet = "Does eating cake, cookies & sweets cause breast cancer?`nThere is no evidence that sugar consumption causesCondition breast cancer - or any other type of cancer. It is true that being overweight can increase your breast cancer risk. Avoiding sugary foods is better for your health!"
$ws.Range("C3").Value = $bigFatSurpriseTweet
$ws.Range("C4").Value = $bigFatSurpriseTweet

$ws.Range("C14").Value = "Reasons why you should include Nuts, Seeds, Olive Oil, Coconut Oil`n" + $bull + "Normalizes Omega 3/6/9 fat ratios`n" + $bull + "Promotes healthy blood flow`n" + $bull + "Reduces body inflammation that causesCondition heart attacks, Alzheimer" + $apos + "s and cancer to name a few`n`nRT and Share"

$ws.Range("C16").Value = "Alcohol, obesity and physical inactivity are all preventable causesCondition of cancer along with tobacco."

# row 14 wraps onto one more line once the text grows, matching Excel's
# automatic re-wrap height for that (wrapped) cell
$ws.Rows.Item(14).RowHeight = 129.6

# leave the cursor where the author's last edit (E17) was
$ws.Range("E17").Select()
